$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Ptprz1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.208684666666667
$ws.Cells.Item(2, 8).Value = 3.626054
$ws.Cells.Item(2, 9).Value = 0.01462795763842055
$ws.Cells.Item(2, 10).Value = 0.01462795763842055
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.138766
$ws.Cells.Item(2, 14).Value = 0.416298
$ws.Cells.Item(2, 15).Value = 0.01356925767068476
$ws.Cells.Item(2, 16).Value = 0.01356925767068476
$ws.Cells.Item(2, 17).Value = 0.1677243364546667
$ws.Cells.Item(2, 18).Value = 1.509519028092
$ws.Cells.Item(2, 19).Value = 0.0001984905263915898
$ws.Cells.Item(2, 20).Value = 0.0001984905263915898

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Ptprz1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.208684666666667
$ws.Cells.Item(3, 8).Value = 3.626054
$ws.Cells.Item(3, 9).Value = 0.01462795763842055
$ws.Cells.Item(3, 10).Value = 0.01462795763842055
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.05416133333333333
$ws.Cells.Item(3, 14).Value = 0.162484
$ws.Cells.Item(3, 15).Value = 0.00529617548814441
$ws.Cells.Item(3, 16).Value = 0.005296175488144411
$ws.Cells.Item(3, 17).Value = 0.06546397312622222
$ws.Cells.Item(3, 18).Value = 0.589175758136
$ws.Cells.Item(3, 19).Value = 0.00007747223068621773
$ws.Cells.Item(3, 20).Value = 0.00007747223068621775

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Ptprz1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.208684666666667
$ws.Cells.Item(4, 8).Value = 3.626054
$ws.Cells.Item(4, 9).Value = 0.01462795763842055
$ws.Cells.Item(4, 10).Value = 0.01462795763842055
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.002366
$ws.Cells.Item(4, 14).Value = 30.007098
$ws.Cells.Item(4, 15).Value = 0.9780831152479456
$ws.Cells.Item(4, 16).Value = 0.9780831152479456
$ws.Cells.Item(4, 17).Value = 12.089706414588
$ws.Cells.Item(4, 18).Value = 108.807357731292
$ws.Cells.Item(4, 19).Value = 0.01430735837670136
$ws.Cells.Item(4, 20).Value = 0.01430735837670136

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Ptprz1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.208684666666667
$ws.Cells.Item(5, 8).Value = 3.626054
$ws.Cells.Item(5, 9).Value = 0.01462795763842055
$ws.Cells.Item(5, 10).Value = 0.01462795763842055
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.03120566666666667
$ws.Cells.Item(5, 14).Value = 0.09361700000000001
$ws.Cells.Item(5, 15).Value = 0.003051451593225274
$ws.Cells.Item(5, 16).Value = 0.003051451593225274
$ws.Cells.Item(5, 17).Value = 0.03771781081311111
$ws.Cells.Item(5, 18).Value = 0.339460297318
$ws.Cells.Item(5, 19).Value = 0.00004463650464139022
$ws.Cells.Item(5, 20).Value = 0.00004463650464139022

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Ptprz1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.178030666666667
$ws.Cells.Item(6, 8).Value = 15.534092
$ws.Cells.Item(6, 9).Value = 0.0626664797952065
$ws.Cells.Item(6, 10).Value = 0.06266647979520648
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.138766
$ws.Cells.Item(6, 14).Value = 0.416298
$ws.Cells.Item(6, 15).Value = 0.01356925767068476
$ws.Cells.Item(6, 16).Value = 0.01356925767068476
$ws.Cells.Item(6, 17).Value = 0.7185346034906668
$ws.Cells.Item(6, 18).Value = 6.466811431416001
$ws.Cells.Item(6, 19).Value = 0.0008503376116559171
$ws.Cells.Item(6, 20).Value = 0.0008503376116559171

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Ptprz1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.178030666666667
$ws.Cells.Item(7, 8).Value = 15.534092
$ws.Cells.Item(7, 9).Value = 0.0626664797952065
$ws.Cells.Item(7, 10).Value = 0.06266647979520648
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.05416133333333333
$ws.Cells.Item(7, 14).Value = 0.162484
$ws.Cells.Item(7, 15).Value = 0.00529617548814441
$ws.Cells.Item(7, 16).Value = 0.005296175488144411
$ws.Cells.Item(7, 17).Value = 0.2804490449475556
$ws.Cells.Item(7, 18).Value = 2.524041404528
$ws.Cells.Item(7, 19).Value = 0.0003318926742196696
$ws.Cells.Item(7, 20).Value = 0.0003318926742196696

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Ncam1"
$ws.Cells.Item(8, 3).Value = "Ptprz1"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.178030666666667
$ws.Cells.Item(8, 8).Value = 15.534092
$ws.Cells.Item(8, 9).Value = 0.0626664797952065
$ws.Cells.Item(8, 10).Value = 0.06266647979520648
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 10.002366
$ws.Cells.Item(8, 14).Value = 30.007098
$ws.Cells.Item(8, 15).Value = 0.9780831152479456
$ws.Cells.Item(8, 16).Value = 0.9780831152479456
$ws.Cells.Item(8, 17).Value = 51.79255788722401
$ws.Cells.Item(8, 18).Value = 466.133020985016
$ws.Cells.Item(8, 19).Value = 0.06129302577971801
$ws.Cells.Item(8, 20).Value = 0.06129302577971799

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Ncam1"
$ws.Cells.Item(9, 3).Value = "Ptprz1"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.178030666666667
$ws.Cells.Item(9, 8).Value = 15.534092
$ws.Cells.Item(9, 9).Value = 0.0626664797952065
$ws.Cells.Item(9, 10).Value = 0.06266647979520648
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.03120566666666667
$ws.Cells.Item(9, 14).Value = 0.09361700000000001
$ws.Cells.Item(9, 15).Value = 0.003051451593225274
$ws.Cells.Item(9, 16).Value = 0.003051451593225274
$ws.Cells.Item(9, 17).Value = 0.1615838989737778
$ws.Cells.Item(9, 18).Value = 1.454255090764
$ws.Cells.Item(9, 19).Value = 0.0001912237296129023
$ws.Cells.Item(9, 20).Value = 0.0001912237296129022

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Ncam1"
$ws.Cells.Item(10, 3).Value = "Ptprz1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 76.16218566666667
$ws.Cells.Item(10, 8).Value = 228.486557
$ws.Cells.Item(10, 9).Value = 0.9217434921665711
$ws.Cells.Item(10, 10).Value = 0.921743492166571
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.138766
$ws.Cells.Item(10, 14).Value = 0.416298
$ws.Cells.Item(10, 15).Value = 0.01356925767068476
$ws.Cells.Item(10, 16).Value = 0.01356925767068476
$ws.Cells.Item(10, 17).Value = 10.56872185622067
$ws.Cells.Item(10, 18).Value = 95.118496705986
$ws.Cells.Item(10, 19).Value = 0.012507374951485
$ws.Cells.Item(10, 20).Value = 0.012507374951485

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ncam1"
$ws.Cells.Item(11, 3).Value = "Ptprz1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 76.16218566666667
$ws.Cells.Item(11, 8).Value = 228.486557
$ws.Cells.Item(11, 9).Value = 0.9217434921665711
$ws.Cells.Item(11, 10).Value = 0.921743492166571
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.05416133333333333
$ws.Cells.Item(11, 14).Value = 0.162484
$ws.Cells.Item(11, 15).Value = 0.00529617548814441
$ws.Cells.Item(11, 16).Value = 0.005296175488144411
$ws.Cells.Item(11, 17).Value = 4.125045525287556
$ws.Cells.Item(11, 18).Value = 37.125409727588
$ws.Cells.Item(11, 19).Value = 0.004881715289569224
$ws.Cells.Item(11, 20).Value = 0.004881715289569224

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ncam1"
$ws.Cells.Item(12, 3).Value = "Ptprz1"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 76.16218566666667
$ws.Cells.Item(12, 8).Value = 228.486557
$ws.Cells.Item(12, 9).Value = 0.9217434921665711
$ws.Cells.Item(12, 10).Value = 0.921743492166571
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 10.002366
$ws.Cells.Item(12, 14).Value = 30.007098
$ws.Cells.Item(12, 15).Value = 0.9780831152479456
$ws.Cells.Item(12, 16).Value = 0.9780831152479456
$ws.Cells.Item(12, 17).Value = 761.8020563979541
$ws.Cells.Item(12, 18).Value = 6856.218507581586
$ws.Cells.Item(12, 19).Value = 0.9015417462778003
$ws.Cells.Item(12, 20).Value = 0.9015417462778001

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ncam1"
$ws.Cells.Item(13, 3).Value = "Ptprz1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 76.16218566666667
$ws.Cells.Item(13, 8).Value = 228.486557
$ws.Cells.Item(13, 9).Value = 0.9217434921665711
$ws.Cells.Item(13, 10).Value = 0.921743492166571
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.03120566666666667
$ws.Cells.Item(13, 14).Value = 0.09361700000000001
$ws.Cells.Item(13, 15).Value = 0.003051451593225274
$ws.Cells.Item(13, 16).Value = 0.003051451593225274
$ws.Cells.Item(13, 17).Value = 2.376691778518778
$ws.Cells.Item(13, 18).Value = 21.390226006669
$ws.Cells.Item(13, 19).Value = 0.002812655647716711
$ws.Cells.Item(13, 20).Value = 0.002812655647716711

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Ncam1"
$ws.Cells.Item(14, 3).Value = "Ptprz1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.07949433333333333
$ws.Cells.Item(14, 8).Value = 0.238483
$ws.Cells.Item(14, 9).Value = 0.0009620703998019471
$ws.Cells.Item(14, 10).Value = 0.000962070399801947
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.138766
$ws.Cells.Item(14, 14).Value = 0.416298
$ws.Cells.Item(14, 15).Value = 0.01356925767068476
$ws.Cells.Item(14, 16).Value = 0.01356925767068476
$ws.Cells.Item(14, 17).Value = 0.01103111065933333
$ws.Cells.Item(14, 18).Value = 0.099279995934
$ws.Cells.Item(14, 19).Value = 0.00001305458115225132
$ws.Cells.Item(14, 20).Value = 0.00001305458115225132

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Ncam1"
$ws.Cells.Item(15, 3).Value = "Ptprz1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.07949433333333333
$ws.Cells.Item(15, 8).Value = 0.238483
$ws.Cells.Item(15, 9).Value = 0.0009620703998019471
$ws.Cells.Item(15, 10).Value = 0.000962070399801947
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.05416133333333333
$ws.Cells.Item(15, 14).Value = 0.162484
$ws.Cells.Item(15, 15).Value = 0.00529617548814441
$ws.Cells.Item(15, 16).Value = 0.005296175488144411
$ws.Cells.Item(15, 17).Value = 0.004305519085777778
$ws.Cells.Item(15, 18).Value = 0.038749671772
$ws.Cells.Item(15, 19).Value = 0.000005095293669300365
$ws.Cells.Item(15, 20).Value = 0.000005095293669300365

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Ncam1"
$ws.Cells.Item(16, 3).Value = "Ptprz1"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.07949433333333333
$ws.Cells.Item(16, 8).Value = 0.238483
$ws.Cells.Item(16, 9).Value = 0.0009620703998019471
$ws.Cells.Item(16, 10).Value = 0.000962070399801947
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 10.002366
$ws.Cells.Item(16, 14).Value = 30.007098
$ws.Cells.Item(16, 15).Value = 0.9780831152479456
$ws.Cells.Item(16, 16).Value = 0.9780831152479456
$ws.Cells.Item(16, 17).Value = 0.795131416926
$ws.Cells.Item(16, 18).Value = 7.156182752334
$ws.Cells.Item(16, 19).Value = 0.0009409848137261249
$ws.Cells.Item(16, 20).Value = 0.0009409848137261247

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Ncam1"
$ws.Cells.Item(17, 3).Value = "Ptprz1"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.07949433333333333
$ws.Cells.Item(17, 8).Value = 0.238483
$ws.Cells.Item(17, 9).Value = 0.0009620703998019471
$ws.Cells.Item(17, 10).Value = 0.000962070399801947
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.03120566666666667
$ws.Cells.Item(17, 14).Value = 0.09361700000000001
$ws.Cells.Item(17, 15).Value = 0.003051451593225274
$ws.Cells.Item(17, 16).Value = 0.003051451593225274
$ws.Cells.Item(17, 17).Value = 0.002480673667888889
$ws.Cells.Item(17, 18).Value = 0.022326063011
$ws.Cells.Item(17, 19).Value = 0.000002935711254270528
$ws.Cells.Item(17, 20).Value = 0.000002935711254270527
